$wb = $excel.ActiveWorkbook

# ALC row 17 (Leve Item ID 38956)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1074.9
$ws.Range("J17").Value = 1074.9
$ws.Range("L17").Value = 3224.7
$ws.Range("N17").Value = -3560.7

# ALC row 18 (Leve Item ID 5471)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1649.2858
$ws.Range("I18").Value = 682.5
$ws.Range("J18").Value = 2938.3333
$ws.Range("K18").Value = 682.5
$ws.Range("L18").Value = 2938.3333
$ws.Range("M18").Value = -398.5
$ws.Range("N18").Value = -3506.3333

# ALC row 19 (Leve Item ID 7015)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1073.4
$ws.Range("I19").Value = 1732.6666
$ws.Range("J19").Value = 790.8570999999999
$ws.Range("K19").Value = 1732.6666
$ws.Range("L19").Value = 790.8570999999999
$ws.Range("M19").Value = -1557.6666
$ws.Range("N19").Value = -1140.8571

# ALC row 113 (Leve Item ID 27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1769.9
$ws.Range("I113").Value = 1744.3334
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1744.3334
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1509.6666
$ws.Range("N113").Value = -8508

# ARM row 2 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1368.5
$ws.Range("I2").Value = 1368.5
$ws.Range("K2").Value = 1368.5
$ws.Range("M2").Value = -1255.5

# ARM row 5 (Leve Item ID 5091)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 150
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 150
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -38
$ws.Range("N5").ClearContents()

# ARM row 116 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1368.5
$ws.Range("I116").Value = 1368.5
$ws.Range("K116").Value = 1368.5
$ws.Range("M116").Value = 925.5

# BSM row 3 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1368.5
$ws.Range("I3").Value = 1368.5
$ws.Range("K3").Value = 1368.5
$ws.Range("M3").Value = -1254.5

# BSM row 4 (Leve Item ID 5091)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -35
$ws.Range("N4").ClearContents()

# BSM row 12 (Leve Item ID 2392)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

# BSM row 17 (Leve Item ID 2393)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# BSM row 29 (Leve Item ID 2318)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 18900
$ws.Range("I29").Value = 27800
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 27800
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -27511
$ws.Range("N29").Value = -10578

# BSM row 134 (Leve Item ID 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1444.4
$ws.Range("I134").Value = 1250
$ws.Range("J134").Value = 2222
$ws.Range("K134").Value = 3750
$ws.Range("L134").Value = 6666
$ws.Range("M134").Value = -1215
$ws.Range("N134").Value = -11736

# CRP row 107 (Leve Item ID 27689)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1882.2
$ws.Range("I107").Value = 2002.75
$ws.Range("K107").Value = 2002.75
$ws.Range("M107").Value = -82.75

# CRP row 122 (Leve Item ID 36196)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1382
$ws.Range("I122").Value = 1227.625
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 3682.875
$ws.Range("L122").Value = 5998.5
$ws.Range("M122").Value = -1232.875
$ws.Range("N122").Value = -10898.5

# CRP row 132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7197.6
$ws.Range("I132").Value = 5994
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 17982
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -15452
$ws.Range("N132").Value = -29060

# CRP row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3479.9333
$ws.Range("I134").Value = 885.7143
$ws.Range("J134").Value = 5749.875
$ws.Range("K134").Value = 2657.1429
$ws.Range("L134").Value = 17249.625
$ws.Range("M134").Value = -122.1428999999998
$ws.Range("N134").Value = -22319.625

# CUL row 2 (Leve Item ID 4847)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 15
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 90
$ws.Range("L2").Value = 90
$ws.Range("M2").Value = 23
$ws.Range("N2").Value = -316

# GSM row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1392
$ws.Range("I122").Value = 1448.7142
$ws.Range("J122").Value = 995
$ws.Range("K122").Value = 4346.142599999999
$ws.Range("L122").Value = 2985
$ws.Range("M122").Value = -1896.142599999999
$ws.Range("N122").Value = -7885

# GSM row 126 (Leve Item ID 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 600
$ws.Range("I126").Value = 600
$ws.Range("K126").Value = 1800
$ws.Range("M126").Value = 670

# GSM row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1484.8889
$ws.Range("I132").Value = 1484.8889
$ws.Range("K132").Value = 4454.6667
$ws.Range("M132").Value = -1924.6667

# GSM row 135 (Leve Item ID 42006)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 92000
$ws.Range("J135").Value = 92000
$ws.Range("L135").Value = 92000
$ws.Range("N135").Value = -102140

# LTW row 96 (Leve Item ID 19735)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# LTW row 100 (Leve Item ID 19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

# LTW row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 14750
$ws.Range("I136").Value = 9999.166999999999
$ws.Range("K136").Value = 29997.501
$ws.Range("M136").Value = -27447.501

# WVR row 2 (Leve Item ID 3307)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 11250
$ws.Range("I2").Value = 11250
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 11250
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -11138
$ws.Range("N2").ClearContents()

# WVR row 33 (Leve Item ID 2734)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

# WVR row 36 (Leve Item ID 2734)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# WVR row 61 (Leve Item ID 2854)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 40000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 40000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 40000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -40584

# WVR row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2833.6667
$ws.Range("I136").Value = 2358.2856
$ws.Range("K136").Value = 7074.8568
$ws.Range("M136").Value = -4524.8568
